# Anonymize test_data: replace the real patient/sample ID in A2
# ("RX232.2022 ") with an anonymous numeric id, and move the active
# selection from D6 to D5 (as recorded in the saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 held the shared string "RX232.2022 " (the sample/patient code).
# Replacing it with a plain number removes that unique string from the
# shared-strings table entirely (count/uniqueCount both drop by one),
# which is exactly the anonymization performed in the commit.
$ws.Range("A2").Value = 1

# The saved cursor/selection position moved up one row.
$ws.Range("D5").Select() | Out-Null
